$d = $word.ActiveDocument

# --- Paragraph: "Par défaut, easyOMR va grouper ..." --------------------
# The paragraph already starts with plain text, so a straight Find &
# Replace of the whole visible sentence (old text == new text) is enough
# to make Word re-flow it into a single run and drop the proofErr marks
# that used to bracket the two "easyOMR" spell-check exceptions.
$p1 = "Par défaut, easyOMR va grouper toutes les cases sur une même ligne horizontale pour former une question. Vous pouvez changer ce groupement grâce au bouton Edition… dans la partie Template d’easyOMR."
$d.Content.Find.Execute($p1, $true, $false, $false, $false, $false, $true, 1, $false, $p1, 2)

# --- Paragraph: "easyOMR ne peut pas travailler ..." ---------------------
# Here the paragraph's very first content is the flagged word "easyOMR",
# so the spellStart/gramStart proofErr marks sit before any run, with
# nothing of the paragraph preceding them. A plain Find & Replace across
# the sentence therefore leaves those two marks behind (nothing inside
# the edited range precedes them). Work around it by temporarily
# inserting a throw-away character right in front of "easyOMR" so the
# proofErr marks become "interior" to the run sequence touched by the
# next replace, then strip that throw-away character back out again.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "easyOMR ne peut pas travailler*") {
        $targetPara = $para
    }
}
if ($targetPara -eq $null) {
    throw "Could not locate the 'easyOMR ne peut pas travailler...' paragraph"
}
$startPos = $targetPara.Range.Start
$d.Range($startPos, $startPos).InsertBefore("Z")

$p2old = "ZeasyOMR ne peut pas travailler directement avec ce document Word, mais pourra travailler avec sa version PDF :"
$d.Content.Find.Execute($p2old, $true, $false, $false, $false, $false, $true, 1, $false, $p2old, 2)
$d.Content.Find.Execute("ZeasyOMR", $true, $false, $false, $false, $false, $true, 1, $false, "easyOMR", 2)

# --- Paragraph: "Sélectionnez Fichier/Enregistrer sous ..." -------------
$p3 = "Sélectionnez Fichier/Enregistrer sous/Double-cliquez Ce PC/Changez Type à PDF (*.pdf) et sauvez dans le répertoire de template du projet, en remplaçant le fichier pdf qui y est présent."
$d.Content.Find.Execute($p3, $true, $false, $false, $false, $false, $true, 1, $false, $p3, 2)
